# Generate Report for Handback
# Update the timestamp cells that record when the handoff/handback XLIFF
# generation happened. These are stored as plain text (not real dates),
# so assign them as strings to avoid Excel re-interpreting them as
# date/time serial numbers.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G), row 2
$wsOverview.Range("G2").Value = "2016-08-18 17:05:23"

# zh-cn sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K), row 2
$wsZhCn.Range("H2").Value = "2016-08-18 17:05:18"
$wsZhCn.Range("K2").Value = "2016-08-18 17:05:35"

# de-de sheet: same two columns, row 2
$wsDeDe.Range("H2").Value = "2016-08-18 17:05:23"
$wsDeDe.Range("K2").Value = "2016-08-18 17:05:42"
